$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "simulation results"

# Reset the view: remove frozen/top-left scroll offset and change selection
$ws.Range("J9").Select()
